$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9813833832740784
$ws.Range("B1").Value = 1.295290946960449
$ws.Range("C1").Value = 2.072923898696899
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 2.021457433700562
